$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 22:19"

# Country name re-ranking (rows whose country changed due to updated case counts)
$ws.Cells.Item(160, 1).Value = "Aruba"
$ws.Cells.Item(161, 1).Value = "Lesoto"
$ws.Cells.Item(162, 1).Value = "Reunion"
$ws.Cells.Item(163, 1).Value = "Crucero"
$ws.Cells.Item(178, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(179, 1).Value = "Eritrea"
$ws.Cells.Item(180, 1).Value = "Camboya"
$ws.Cells.Item(181, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(182, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$updates = @(
    @{Row=4; B=5343688; C=37731; D=2787473; E=2387554; F=0; G=916; H=168661},
    @{Row=8; B=568919; C=2810; D=432029; E=125880; F=0; G=259; H=11010},
    @{Row=27; B=120633; C=212; D=107023; E=4606; F=0; G=13; H=9004},
    @{Row=69; B=25057; C=549; D=8189; E=16605; F=0; G=8; H=263},
    @{Row=89; B=8423; C=63; D=7713; E=660; F=0; G=1; H=50},
    @{Row=91; B=8116; C=98; D=7060; E=1006; F=0; G=0; H=50},
    @{Row=102; B=6177; C=235; D=3804; E=2157; F=0; G=2; H=216},
    @{Row=109; B=4752; C=38; D=2529; E=2071; F=0; G=0; H=152},
    @{Row=110; B=4652; C=7; D=1728; E=2863; F=0; G=0; H=61},
    @{Row=115; B=3525; C=115; D=1910; E=1552; F=0; G=0; H=63},
    @{Row=127; B=2477; C=5; D=1175; E=1255; F=0; G=0; H=47},
    @{Row=160; B=798; C=81; D=114; E=681; F=0; G=0; H=3},
    @{Row=161; B=781; C=0; D=175; E=582; F=0; G=0; H=24},
    @{Row=162; B=734; C=32; D=631; E=98; F=0; G=0; H=5},
    @{Row=163; B=712; C=0; D=651; E=48; F=0; G=0; H=13},
    @{Row=178; B=287; C=73; D=78; E=206; F=0; G=0; H=3},
    @{Row=179; B=285; C=0; D=248; E=37; F=0; G=0; H=0},
    @{Row=180; B=268; C=2; D=220; E=48; F=0; G=0; H=0},
    @{Row=181; B=224; C=8; D=39; E=183; F=0; G=0; H=2},
    @{Row=182; B=219; C=14; D=102; E=100; F=0; G=0; H=17},
    @{Row=189; B=141; C=3; D=114; E=23; F=0; G=0; H=4},
    @{Row=194; B=90; C=1; D=87; E=2; F=0; G=0; H=1},
    @{Row=213; B=13; C=0; D=12; E=0; F=0; G=0; H=1},
    @{Row=214; B=13; C=0; D=13; E=0; F=0; G=0; H=0}
)

foreach ($r in $updates) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
